$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.142924189567566
$ws.Range("B1").Value = 2.211016416549683
$ws.Range("D1").Value = 2.229090929031372
$ws.Range("E1").Value = 1.077126383781433
